# Update market-data derived columns (H-N) across sheets per scheduled runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 846.4167
$ws.Range("I80").Value = 636.2857
$ws.Range("J80").Value = 1140.6
$ws.Range("K80").Value = 1908.8571
$ws.Range("L80").Value = 3421.8
$ws.Range("M80").Value = -910.8571000000002
$ws.Range("N80").Value = -5417.799999999999
$ws.Range("H83").Value = 846.4167
$ws.Range("I83").Value = 636.2857
$ws.Range("J83").Value = 1140.6
$ws.Range("K83").Value = 5726.571300000001
$ws.Range("L83").Value = 10265.4
$ws.Range("M83").Value = -734.5713000000005
$ws.Range("N83").Value = -20249.4
$ws.Range("H113").Value = 4000
$ws.Range("J113").Value = 4000
$ws.Range("L113").Value = 4000
$ws.Range("N113").Value = -10508

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7475.5
$ws.Range("I61").Value = 7475.5
$ws.Range("K61").Value = 7475.5
$ws.Range("M61").Value = -7263.5
$ws.Range("H88").Value = 1539.3334
$ws.Range("I88").Value = 1392.3334
$ws.Range("J88").Value = 1833.3334
$ws.Range("K88").Value = 1392.3334
$ws.Range("L88").Value = 1833.3334
$ws.Range("M88").Value = -986.3334
$ws.Range("N88").Value = -2645.3334
$ws.Range("H91").Value = 1539.3334
$ws.Range("I91").Value = 1392.3334
$ws.Range("J91").Value = 1833.3334
$ws.Range("K91").Value = 1392.3334
$ws.Range("L91").Value = 1833.3334
$ws.Range("M91").Value = 11.66660000000002
$ws.Range("N91").Value = -4641.3334
$ws.Range("H132").Value = 8999.286
$ws.Range("I132").Value = 1665
$ws.Range("J132").Value = 14500
$ws.Range("K132").Value = 4995
$ws.Range("L132").Value = 43500
$ws.Range("M132").Value = -2465
$ws.Range("N132").Value = -48560
$ws.Range("H136").Value = 7475.5
$ws.Range("I136").Value = 7475.5
$ws.Range("K136").Value = 22426.5
$ws.Range("M136").Value = -19876.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 4006
$ws.Range("I15").Value = 4006
$ws.Range("K15").Value = 4006
$ws.Range("M15").Value = -3779
$ws.Range("H86").Value = 775
$ws.Range("I86").Value = 775
$ws.Range("K86").Value = 775
$ws.Range("M86").Value = 348
$ws.Range("H89").Value = 775
$ws.Range("I89").Value = 775
$ws.Range("K89").Value = 3875
$ws.Range("M89").Value = 1741

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 457.8
$ws.Range("I2").Value = 49.5
$ws.Range("J2").Value = 730
$ws.Range("K2").Value = 49.5
$ws.Range("L2").Value = 730
$ws.Range("M2").Value = 63.5
$ws.Range("N2").Value = -956
$ws.Range("H11").Value = 95
$ws.Range("I11").Value = 100
$ws.Range("J11").Value = 90
$ws.Range("K11").Value = 100
$ws.Range("L11").Value = 90
$ws.Range("M11").Value = 40
$ws.Range("N11").Value = -370
$ws.Range("H13").Value = 994.75
$ws.Range("J13").Value = 994.75
$ws.Range("L13").Value = 994.75
$ws.Range("N13").Value = -1272.75
$ws.Range("H14").Value = 748.75
$ws.Range("J14").Value = 765
$ws.Range("L14").Value = 765
$ws.Range("N14").Value = -1105
$ws.Range("H15").Value = 1050
$ws.Range("J15").Value = 866.6667
$ws.Range("L15").Value = 866.6667
$ws.Range("N15").Value = -1206.6667
$ws.Range("H36").Value = 3160
$ws.Range("I36").Value = 2240
$ws.Range("K36").Value = 2240
$ws.Range("M36").Value = -1852
$ws.Range("H40").Value = 3160
$ws.Range("I40").Value = 2240
$ws.Range("K40").Value = 2240
$ws.Range("M40").Value = -2080
$ws.Range("H42").Value = 220
$ws.Range("I42").Value = 220
$ws.Range("K42").Value = 220
$ws.Range("M42").Value = 373
$ws.Range("H43").Value = 16727
$ws.Range("J43").Value = 16727
$ws.Range("L43").Value = 16727
$ws.Range("N43").Value = -17095
$ws.Range("M56").ClearContents()
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("H82").Value = 65000
$ws.Range("J82").Value = 65000
$ws.Range("L82").Value = 65000
$ws.Range("N82").Value = -65722
$ws.Range("H85").Value = 65000
$ws.Range("J85").Value = 65000
$ws.Range("L85").Value = 65000
$ws.Range("N85").Value = -67496
$ws.Range("H101").Value = 16727
$ws.Range("J101").Value = 16727
$ws.Range("L101").Value = 16727
$ws.Range("N101").Value = -23217
$ws.Range("H132").Value = 4430.857
$ws.Range("I132").Value = 1003.2
$ws.Range("K132").Value = 3009.6
$ws.Range("M132").Value = -479.6000000000004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 33.875
$ws.Range("I12").Value = 38.5
$ws.Range("J12").Value = 29.25
$ws.Range("K12").Value = 115.5
$ws.Range("L12").Value = 87.75
$ws.Range("M12").Value = 57.5
$ws.Range("N12").Value = -433.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 35000
$ws.Range("J15").Value = 35000
$ws.Range("L15").Value = 35000
$ws.Range("N15").Value = -35576
$ws.Range("H47").Value = 35000
$ws.Range("J47").Value = 35000
$ws.Range("L47").Value = 35000
$ws.Range("N47").Value = -36136
$ws.Range("H81").Value = 35000
$ws.Range("J81").Value = 35000
$ws.Range("L81").Value = 35000
$ws.Range("N81").Value = -36996
$ws.Range("H84").Value = 35000
$ws.Range("J84").Value = 35000
$ws.Range("L84").Value = 105000
$ws.Range("N84").Value = -114984
$ws.Range("N97").ClearContents()
$ws.Range("H97").Value = 600
$ws.Range("I97").Value = 600
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 600
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -104

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M4").ClearContents()
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("H22").Value = 3150.25
$ws.Range("J22").Value = 5500
$ws.Range("L22").Value = 5500
$ws.Range("N22").Value = -6090
$ws.Range("H27").Value = 3150.25
$ws.Range("J27").Value = 5500
$ws.Range("L27").Value = 5500
$ws.Range("N27").Value = -5714
$ws.Range("M28").ClearContents()
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("H46").Value = 6170.5884
$ws.Range("J46").Value = 6062.5
$ws.Range("L46").Value = 6062.5
$ws.Range("N46").Value = -6438.5
$ws.Range("H132").Value = 15657.2
$ws.Range("I132").Value = 10734.6
$ws.Range("J132").Value = 20579.8
$ws.Range("K132").Value = 32203.8
$ws.Range("L132").Value = 61739.39999999999
$ws.Range("M132").Value = -29673.8
$ws.Range("N132").Value = -66799.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M70").ClearContents()
$ws.Range("H70").Value = 25000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H73").Value = 25000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("H80").Value = 70000
$ws.Range("J80").Value = 70000
$ws.Range("L80").Value = 70000
$ws.Range("N80").Value = -71996
$ws.Range("H83").Value = 70000
$ws.Range("J83").Value = 70000
$ws.Range("L83").Value = 210000
$ws.Range("N83").Value = -219984
$ws.Range("M94").ClearContents()
$ws.Range("H94").Value = 29000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 29000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 29000
$ws.Range("N94").Value = -30802
$ws.Range("H103").Value = 19866.666
$ws.Range("J103").Value = 19866.666
$ws.Range("L103").Value = 19866.666
$ws.Range("N103").Value = -22210.666
$ws.Range("H132").Value = 11686.667
$ws.Range("I132").Value = 9795.75
$ws.Range("J132").Value = 13199.4
$ws.Range("K132").Value = 29387.25
$ws.Range("L132").Value = 39598.2
$ws.Range("M132").Value = -26857.25
$ws.Range("N132").Value = -44658.2
